# IRPTestStudentMapping.xlsx - add a "CAT" column to flag CAT-style test packages.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at C (shifts former C -> D, D -> E) and give it the
# same width as column B, matching what Excel does when inserting next to
# an existing column.
$ws.Columns("C:C").Insert()
$ws.Columns("C:C").ColumnWidth = $ws.Columns("B:B").ColumnWidth

# Header for the new column.
$ws.Range("C2").Value = "CAT"

# Flag the rows whose Test name (column A) refers to a CAT-style package.
$catRows = @(6,7,8,9,10,11,18,19,20,30,31,32,33,34,35)
foreach ($r in $catRows) {
    $ws.Cells.Item($r, 3).Value = $true
}
